# Update the workbook with newly computed TPM-based values.
# Only columns M through T (Receptor/Edge derived metrics) change, for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.676671000000001
$ws.Range("N2").Value = 26.030013
$ws.Range("O2").Value = 0.1325240072999665
$ws.Range("P2").Value = 0.1325240072999665
$ws.Range("Q2").Value = 1.891144073370667
$ws.Range("R2").Value = 17.020296660336
$ws.Range("S2").Value = 0.1325240072999665
$ws.Range("T2").Value = 0.1325240072999665

# Row 3
$ws.Range("M3").Value = 37.74750533333334
$ws.Range("N3").Value = 113.242516
$ws.Range("O3").Value = 0.5765403197090441
$ws.Range("P3").Value = 0.576540319709044
$ws.Range("Q3").Value = 8.227345602439112
$ws.Range("R3").Value = 74.04611042195201
$ws.Range("S3").Value = 0.5765403197090441
$ws.Range("T3").Value = 0.576540319709044

# Row 4
$ws.Range("M4").Value = 19.04827033333333
$ws.Range("N4").Value = 57.144811
$ws.Range("O4").Value = 0.2909356729909895
$ws.Range("P4").Value = 0.2909356729909895
$ws.Range("Q4").Value = 4.151710206465778
$ws.Range("R4").Value = 37.36539185819201
$ws.Range("S4").Value = 0.2909356729909895
$ws.Range("T4").Value = 0.2909356729909895
